$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "stim-info_gratings-v1"

# Add new row 8 label first (to match shared-string insertion order)
$ws.Range("A8").Value = "param7"

# Update existing param labels (column B) to reflect new ordering
$ws.Range("B2").Value = "total_dur_sec"
$ws.Range("B5").Value = "velocity"
$ws.Range("B6").Value = "direction_deg"
$ws.Range("B7").Value = "michelson_contrast"
$ws.Range("B8").Value = "mean_luminance"

# Finish row 8
$ws.Range("C8").Value = "float"

# Update selection to match the recorded view state
$ws.Range("B5").Select()
